# Daily attendance processing - 2025-10-17 18:24:22
# Swap the order of "System" and the recording user's email in the
# "Recorded By" column (G) for the specific rows that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where "System, dnasr281@gmail.com" becomes "dnasr281@gmail.com, System"
$dnasrRows = @(3,6,10,11,12,13,14,15,17,18,19,30,33,37,38,39,40,41,42,44,45,46,57,60,64,65,66,67,68,69,71,72,73,86,87,88,89,93,95,96,97,99,112,113,114,115,119,121,122,123,125,138,139,140,141,145,147,148,149,151)

foreach ($row in $dnasrRows) {
    $cell = $ws.Range("G$row")
    $cell.Value2 = "dnasr281@gmail.com, System"
}

# Rows where "backup@backdoor.com, System" becomes "System, backup@backdoor.com"
$backdoorRows = @(5,32,59,85,111,137)

foreach ($row in $backdoorRows) {
    $cell = $ws.Range("G$row")
    $cell.Value2 = "System, backup@backdoor.com"
}
